$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Sunday hours for the week of row 9
$ws.Range("H9").Value = 5.25

# Update the active selection to match the saved workbook state
$ws.Range("P10").Select()
